$d = $word.ActiveDocument

# Step 1: remove the _GoBack bookmark that currently sits after the
# "El alcance..." paragraph.
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# Step 2: append the new paragraphs (blank spacer, heading, and the
# five acronym/definition lines) at the end of the document body, right
# after the "Ahora, tomando..." paragraph. We build the exact OOXML for
# these paragraphs (preserving run/paragraph formatting, proofErr spans
# and the relocated _GoBack bookmark) and insert it with InsertXML so
# the structure matches precisely.
$rng = $d.Content
$rng.Collapse(0)

$xmlFragment = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="center"/>
        <w:rPr>
          <w:sz w:val="28"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:sz w:val="28"/>
        </w:rPr>
        <w:t>Definiciones, acrónimos y abreviaciones</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>SegV</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>: Significa “Seguridad en Viviendas”.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t xml:space="preserve">- </w:t>
      </w:r>
      <w:proofErr w:type="spellStart"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>SegE</w:t>
      </w:r>
      <w:proofErr w:type="spellEnd"/>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>: Significa “Seguridad en Empresas”.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>- UMSNH: “Universidad Michoacana de San Nicolás de Hidalgo”.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>- Mb: “Megabytes”.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:tabs>
          <w:tab w:val="left" w:pos="3252"/>
        </w:tabs>
        <w:jc w:val="both"/>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:sz w:val="24"/>
        </w:rPr>
        <w:t>- PF: “Policía federal”</w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p>
'@

$rng.InsertXML($xmlFragment)
